$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '73.841.82'
$ws.Range("E2").Value = '  +7.25%  '
$ws.Range("D3").Value = '2.620.03'
$ws.Range("E3").Value = '  +7.45%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '184.38'
$ws.Range("E5").Value = '  +13.78%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '580.86'
$ws.Range("E6").Value = '  +3.72%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").Value = '  -0.07%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.533'
$ws.Range("E8").Value = '  +4.25%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.198'
$ws.Range("E9").Value = '  +17.72%  '
$ws.Range("D10").Value = '2.617.79'
$ws.Range("E10").Value = '  +7.42%  '
$ws.Range("E11").Value = '  +0.13%  '
$ws.Range("E12").Value = '  +8.32%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.75'
$ws.Range("E13").Value = '  +4.31%  '
$ws.Range("B14").Value = 'ShibaInu'
$ws.Range("C14").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000188'
$ws.Range("E14").Value = '  +6.50%  '
$ws.Range("B15").Value = 'WrappedBTC'
$ws.Range("C15").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D15").Value = '73.510.78'
$ws.Range("E15").Value = '  +6.98%  '
$ws.Range("D16").Value = '3.076.46'
$ws.Range("E16").Value = '  +6.59%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '26.07'
$ws.Range("E17").Value = '  +11.94%  '
$ws.Range("D18").Value = '2.622.83'
$ws.Range("E18").Value = '  +7.48%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '9.09'
$ws.Range("E19").Value = '  +31.10%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.84'
$ws.Range("E20").Value = '  +12.17%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '372.31'
$ws.Range("E21").Value = '  +9.71%  '
$ws.Range("E22").Value = '  +16.56%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.08'
$ws.Range("E23").Value = '  +5.85%  '
$ws.Range("E24").Value = '  +0.07%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '69.65'
$ws.Range("E25").Value = '  +3.82%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '4.14'
$ws.Range("E26").Value = '  +11.10%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.32'
$ws.Range("E27").Value = '  +13.52%  '
$ws.Range("E28").Value = '  +6.51%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.01'
$ws.Range("E29").Value = '  +1.26%  '
$ws.Range("E30").Value = '  +13.73%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '515.37'
$ws.Range("E31").Value = '  +20.16%  '
$ws.Range("E32").Value = '  +18.68%  '
$ws.Range("E33").Value = '  +6.25%  '
$ws.Range("E35").Value = '  +0.00%  '
$ws.Range("E36").Value = '  +12.79%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '161.36'
$ws.Range("E37").Value = '  +1.51%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '19.18'
$ws.Range("E38").Value = '  +6.47%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '19.27'
$ws.Range("E39").Value = '  +1.44%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.87'
$ws.Range("E41").Value = '  +11.87%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.66'
$ws.Range("E42").Value = '  +10.13%  '
$ws.Range("E43").Value = '  +8.41%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '160.24'
$ws.Range("E44").Value = '  +23.04%  '
$ws.Range("E45").Value = '  +9.22%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0864'
$ws.Range("E46").Value = '  +20.15%  '
$ws.Range("E47").Value = '  +13.90%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '38.58'
$ws.Range("E48").Value = '  +3.18%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '3.62'
$ws.Range("E50").Value = '  +9.66%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '20.35'
$ws.Range("E51").Value = '  +20.55%  '
